$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 194, pushing existing rows 194:283 down to 195:284.
$ws.Rows.Item(194).Insert()

# Populate the newly inserted row 194 with the new data record.
$ws.Cells.Item(194, 1).Value  = 3
$ws.Cells.Item(194, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(194, 3).Value  = "Coquimbo"
$ws.Cells.Item(194, 4).Value  = 45229
$ws.Cells.Item(194, 5).Value  = 5
$ws.Cells.Item(194, 6).Value  = 100112030
$ws.Cells.Item(194, 7).Value  = "Poroto granado"
$ws.Cells.Item(194, 8).Value  = "Sin especificar"
$ws.Cells.Item(194, 9).Value  = "Primera"
$ws.Cells.Item(194, 10).Value = 38
$ws.Cells.Item(194, 11).Value = 35000
$ws.Cells.Item(194, 12).Value = 35000
$ws.Cells.Item(194, 13).Value = 35000
$ws.Cells.Item(194, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(194, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(194, 16).Value = 1400
$ws.Cells.Item(194, 17).Value = 25
$ws.Cells.Item(194, 18).Value = "Hortaliza"
